$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Insert a new paragraph "#First Readme formatted" and an
# empty paragraph before the existing "Hello" paragraph.
# (a throw-away placeholder "ZZZ" is typed into what will become the
#  empty paragraph and then deleted, so it serializes as a truly
#  empty <w:p/> rather than leaving a stray empty run behind)
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$insertAtStart = $d.Range($firstPara.Range.Start, $firstPara.Range.Start)
$insertAtStart.InsertBefore("#First Readme formatted`rZZZ`r")

$blankPara = $d.Paragraphs(2)
$placeholder = $d.Range($blankPara.Range.Start, $blankPara.Range.Start + 3)
$placeholder.Delete()

# ------------------------------------------------------------------
# Step 2: Append a new paragraph "How are you doing?" at the very end
# of the document (right after the hidden _GoBack bookmark, which
# currently sits at the end of the "Hello" paragraph). Because the
# text is typed exactly at Content.End, the bookmark naturally ends
# up attached to the end of this new, final paragraph once split off.
# ------------------------------------------------------------------
$endOfDoc = $d.Range($d.Content.End, $d.Content.End)
$endOfDoc.InsertAfter("How are you doing?")
$howAreYouStart = $endOfDoc.Start
$splitBeforeHow = $d.Range($howAreYouStart, $howAreYouStart)
$splitBeforeHow.InsertParagraphBefore()

# ------------------------------------------------------------------
# Step 3: Turn "Hello" into two runs - "Hello" and " Poonam  " - by
# temporarily splitting it into two paragraphs, typing the new text
# into the second one, then deleting the paragraph mark that joins
# them back together. This preserves the run boundary between the
# original "Hello" run and the newly typed " Poonam  " run instead of
# them being coalesced into a single run.
# ------------------------------------------------------------------
$helloPara = $d.Paragraphs(3)
$helloRange = $helloPara.Range
$beforeMark = $d.Range($helloRange.End - 1, $helloRange.End - 1)
$beforeMark.InsertParagraphAfter()

$poonamPara = $d.Paragraphs(4)
$poonamPara.Range.InsertBefore(" Poonam  ")

$helloPara2 = $d.Paragraphs(3)
$paraMark = $d.Range($helloPara2.Range.End - 1, $helloPara2.Range.End)
$paraMark.Delete()
